$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hotel Data")

# Row 2: Holiday Inn NAIROBI TWO RIVERS MALL by IHG
$ws.Range("C2").Value = "₹ 20,452"
$ws.Range("D2").Value = "₹ 130,890"

# Row 3: JW Marriott Hotel Nairobi
$ws.Range("C3").Value = "₹ 204,675"
$ws.Range("D3").Value = "₹ 1,292,111"

# Row 4: Yaya Hotel & Apartments
$ws.Range("C4").Value = "₹ 19,783"
$ws.Range("D4").Value = "₹ 95,383"
